# Append the new "24/06/2024" data block (rows 1100-1132) to Sheet1,
# which holds the FII/DII buy-sell watchlist data. This mirrors the
# existing per-day block layout already present in rows 1-1099.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1100").Value = 'Buying Opportunity'
$ws.Range("B1100").Value = 'support Zone'
$ws.Range("C1100").Value = 'long buildup'
$ws.Range("D1100").Value = 'Short buildup'
$ws.Range("E1100").Value = 'FII ENTERING'

$ws.Range("A1101").Value = 'AFFLE'
$ws.Range("B1101").Value = '3IINFOLTD'
$ws.Range("C1101").Value = 'NAUKRI'
$ws.Range("D1101").Value = 'ABBOTINDIA'
$ws.Range("E1101").Value = 'ASTRAL'
$ws.Range("F1101").Value = 1316.8
$ws.Range("G1101").Value = 36.95
$ws.Range("H1101").Value = 6594.1
$ws.Range("I1101").Value = 26834.45
$ws.Range("J1101").Value = 2384.4

$ws.Range("A1102").Value = 'AJMERA'
$ws.Range("B1102").Value = 'EIDPARRY'
$ws.Range("D1102").Value = 'ACC'
$ws.Range("E1102").Value = 'AUBANK'
$ws.Range("F1102").Value = 718.4
$ws.Range("G1102").Value = 772.8
$ws.Range("I1102").Value = 2588.25
$ws.Range("J1102").Value = 679.5

$ws.Range("A1103").Value = 'ALPA'
$ws.Range("B1103").Value = 'GIPCL'
$ws.Range("D1103").Value = 'BPCL'
$ws.Range("E1103").Value = 'CUMMINSIND'
$ws.Range("F1103").Value = 89.97
$ws.Range("G1103").Value = 252.96
$ws.Range("I1103").Value = 305.25
$ws.Range("J1103").Value = 4042.5

$ws.Range("A1104").Value = 'ASAHIINDIA'
$ws.Range("B1104").Value = 'GRAPHITE'
$ws.Range("D1104").Value = 'COALINDIA'
$ws.Range("E1104").Value = 'NAUKRI'
$ws.Range("F1104").Value = 708.85
$ws.Range("G1104").Value = 584.8
$ws.Range("I1104").Value = 473.7
$ws.Range("J1104").Value = 6594.1

$ws.Range("A1105").Value = 'ASTRAL'
$ws.Range("B1105").Value = 'HEALTHY'
$ws.Range("D1105").Value = 'IRCTC'
$ws.Range("E1105").Value = 'PAGEIND'
$ws.Range("F1105").Value = 2384.4
$ws.Range("G1105").Value = 12.63
$ws.Range("I1105").Value = 1010.25
$ws.Range("J1105").Value = 40307.45

$ws.Range("A1106").Value = 'AVALON'
$ws.Range("B1106").Value = 'IFCI'
$ws.Range("F1106").Value = 550.95
$ws.Range("G1106").Value = 60.97

$ws.Range("A1107").Value = 'BAJAJHCARE'
$ws.Range("B1107").Value = 'IRISDOREME'
$ws.Range("F1107").Value = 343.7
$ws.Range("G1107").Value = 74.48

$ws.Range("A1108").Value = 'BBTC'
$ws.Range("B1108").Value = 'KARURVYSYA'
$ws.Range("F1108").Value = 2014.4
$ws.Range("G1108").Value = 209.77

$ws.Range("A1109").Value = 'BOMDYEING'
$ws.Range("B1109").Value = 'MAGNUM'
$ws.Range("F1109").Value = 198.27
$ws.Range("G1109").Value = 52.31

$ws.Range("A1110").Value = 'CGPOWER'
$ws.Range("B1110").Value = 'MANAKSIA'
$ws.Range("F1110").Value = 701.75
$ws.Range("G1110").Value = 97.91

$ws.Range("A1111").Value = 'CUMMINSIND'
$ws.Range("B1111").Value = 'RAMAPHO'
$ws.Range("F1111").Value = 4042.5
$ws.Range("G1111").Value = 188.23

$ws.Range("A1112").Value = 'DCXINDIA'
$ws.Range("B1112").Value = 'RKEC'
$ws.Range("F1112").Value = 372.75
$ws.Range("G1112").Value = 106.07

$ws.Range("A1113").Value = 'DVL'
$ws.Range("F1113").Value = 378.7

$ws.Range("A1114").Value = 'ENIL'
$ws.Range("F1114").Value = 255.59

$ws.Range("A1115").Value = 'GANECOS'
$ws.Range("F1115").Value = 1363

$ws.Range("A1116").Value = 'GOCLCORP'
$ws.Range("F1116").Value = 477.25

$ws.Range("A1117").Value = 'GPPL'
$ws.Range("F1117").Value = 213.78

$ws.Range("A1118").Value = 'GULFOILLUB'
$ws.Range("F1118").Value = 1037.4

$ws.Range("A1119").Value = 'HINDWAREAP'
$ws.Range("F1119").Value = 439.65

$ws.Range("A1120").Value = 'JINDALPOLY'
$ws.Range("F1120").Value = 690.3

$ws.Range("A1121").Value = 'KICL'
$ws.Range("F1121").Value = 5926.95

$ws.Range("A1122").Value = 'NAUKRI'
$ws.Range("F1122").Value = 6594.1

$ws.Range("A1123").Value = 'NEULANDLAB'
$ws.Range("F1123").Value = 6704.85

$ws.Range("A1124").Value = 'NUCLEUS'
$ws.Range("F1124").Value = 1437.15

$ws.Range("A1125").Value = 'OMAXE'
$ws.Range("F1125").Value = 119

$ws.Range("A1126").Value = 'PAGEIND'
$ws.Range("F1126").Value = 40307.45

$ws.Range("A1127").Value = 'PCBL'
$ws.Range("F1127").Value = 274.26

$ws.Range("A1128").Value = 'PDMJEPAPER'
$ws.Range("F1128").Value = 112.05

$ws.Range("A1129").Value = 'PIXTRANS'
$ws.Range("F1129").Value = 1353.45

$ws.Range("A1130").Value = 'PYRAMID'
$ws.Range("F1130").Value = 169.16

$ws.Range("A1131").Value = 'ROTO'
$ws.Range("F1131").Value = 528.75

$ws.Range("A1132").Value = '24/06/2024'
